# Update "想去人数" (interested-count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): row 5 corresponds to the event whose
# interested-count grew from 4314 to 4450.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 4450

# Sheet "演出" (Performances): row 2 corresponds to the event whose
# interested-count grew from 124 to 125.
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 125

# Sheet "全部类型" (All types, aggregate of every other sheet) mirrors the
# same two events in different rows: row 3 (the 演出 event) and row 9
# (the 展览 event).
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 125
$wsAll.Range("F9").Value = 4450
